$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row 1: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2310"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2404"
        }
    }
}

# --- Turn the range into an Excel Table ---
$range = $ws.Range("A1:U85")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- Freeze the header row ---
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Application.ActiveWindow.FreezePanes = $true
